$wb = $excel.ActiveWorkbook

# ALC row 68
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H68").Value = 0
$ws.Range("J68").Value = 0
$ws.Range("L68").Value = 0
$ws.Range("N68").Value = ""

# ALC row 71
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H71").Value = 0
$ws.Range("J71").Value = 0
$ws.Range("L71").Value = 0
$ws.Range("N71").Value = ""

# ALC row 80
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H80").Value = 8215.462
$ws.Range("I80").Value = 560.7778
$ws.Range("K80").Value = 1682.3334
$ws.Range("M80").Value = -684.3334

# ALC row 83
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H83").Value = 8215.462
$ws.Range("I83").Value = 560.7778
$ws.Range("K83").Value = 5047.000199999999
$ws.Range("M83").Value = -55.0001999999995

# ALC row 86
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H86").Value = 5545.727
$ws.Range("I86").Value = 5545.727
$ws.Range("J86").Value = 0
$ws.Range("K86").Value = 5545.727
$ws.Range("L86").Value = 0
$ws.Range("M86").Value = -4422.727
$ws.Range("N86").Value = ""

# ALC row 88
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H88").Value = 21999
$ws.Range("I88").Value = 2997
$ws.Range("K88").Value = 2997
$ws.Range("M88").Value = -2591

# ALC row 89
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H89").Value = 5545.727
$ws.Range("I89").Value = 5545.727
$ws.Range("J89").Value = 0
$ws.Range("K89").Value = 27728.635
$ws.Range("L89").Value = 0
$ws.Range("M89").Value = -22112.635
$ws.Range("N89").Value = ""

# ALC row 91
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H91").Value = 21999
$ws.Range("I91").Value = 2997
$ws.Range("K91").Value = 2997
$ws.Range("M91").Value = -1593

# ALC row 92
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H92").Value = 35354480
$ws.Range("I92").Value = 9260153
$ws.Range("J92").Value = 66667668
$ws.Range("K92").Value = 9260153
$ws.Range("L92").Value = 66667668
$ws.Range("M92").Value = -9258905
$ws.Range("N92").Value = -66670164

# ALC row 125
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H125").Value = 3118.3333
$ws.Range("I125").Value = 3383.125
$ws.Range("J125").Value = 1000
$ws.Range("K125").Value = 30448.125
$ws.Range("L125").Value = 9000
$ws.Range("M125").Value = -27988.125
$ws.Range("N125").Value = -13920

# ALC row 132
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H132").Value = 4506411
$ws.Range("I132").Value = 1529.7301
$ws.Range("J132").Value = 30307094
$ws.Range("K132").Value = 4589.1903
$ws.Range("L132").Value = 90921282
$ws.Range("M132").Value = -2059.1903
$ws.Range("N132").Value = -90926342

# ALC row 137
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H137").Value = 1321.7142
$ws.Range("I137").Value = 1352.4546
$ws.Range("J137").Value = 1269.6923
$ws.Range("K137").Value = 4057.3638
$ws.Range("L137").Value = 3809.0769
$ws.Range("M137").Value = -1507.3638
$ws.Range("N137").Value = -8909.0769

# ARM row 32
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 11043.24
$ws.Range("I32").Value = 9413.474
$ws.Range("J32").Value = 16204.167
$ws.Range("K32").Value = 9413.474
$ws.Range("L32").Value = 16204.167
$ws.Range("M32").Value = -9126.474
$ws.Range("N32").Value = -16778.167

# ARM row 74
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 12501341
$ws.Range("I74").Value = 974.86957
$ws.Range("J74").Value = 29413600
$ws.Range("K74").Value = 974.86957
$ws.Range("L74").Value = 29413600
$ws.Range("M74").Value = -100.86957
$ws.Range("N74").Value = -29415348

# ARM row 77
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H77").Value = 12501341
$ws.Range("I77").Value = 974.86957
$ws.Range("J77").Value = 29413600
$ws.Range("K77").Value = 4874.34785
$ws.Range("L77").Value = 147068000
$ws.Range("M77").Value = -506.3478500000001
$ws.Range("N77").Value = -147076736

# ARM row 88
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H88").Value = 2609.5
$ws.Range("I88").Value = 0
$ws.Range("K88").Value = 0
$ws.Range("M88").Value = ""

# ARM row 91
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H91").Value = 2609.5
$ws.Range("I91").Value = 0
$ws.Range("K91").Value = 0
$ws.Range("M91").Value = ""

# ARM row 110
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H110").Value = 7175.727
$ws.Range("I110").Value = 8617.790999999999
$ws.Range("J110").Value = 2008.3334
$ws.Range("K110").Value = 8617.790999999999
$ws.Range("L110").Value = 2008.3334
$ws.Range("M110").Value = -6572.790999999999
$ws.Range("N110").Value = -6098.3334

# ARM row 122
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H122").Value = 4284339.5
$ws.Range("I122").Value = 6424009
$ws.Range("K122").Value = 19272027
$ws.Range("M122").Value = -19269577

# BSM row 20
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 9452.718999999999
$ws.Range("I20").Value = 1676.2727
$ws.Range("J20").Value = 26560.9
$ws.Range("K20").Value = 1676.2727
$ws.Range("L20").Value = 26560.9
$ws.Range("M20").Value = -1429.2727
$ws.Range("N20").Value = -27054.9

# CRP row 58
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 7362900.5
$ws.Range("I58").Value = 11906703
$ws.Range("J58").Value = 1001577.5
$ws.Range("K58").Value = 11906703
$ws.Range("L58").Value = 1001577.5
$ws.Range("M58").Value = -11906500
$ws.Range("N58").Value = -1001983.5

# CRP row 97
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H97").Value = 19700
$ws.Range("J97").Value = 19700
$ws.Range("L97").Value = 19700
$ws.Range("N97").Value = -21682

# CRP row 132
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 5884016.5
$ws.Range("I132").Value = 8334735
$ws.Range("K132").Value = 25004205
$ws.Range("M132").Value = -25001675

# CRP row 134
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H134").Value = 11078312
$ws.Range("I134").Value = 14496523
$ws.Range("J134").Value = 1250954.9
$ws.Range("K134").Value = 43489569
$ws.Range("L134").Value = 3752864.7
$ws.Range("M134").Value = -43487034
$ws.Range("N134").Value = -3757934.7

# CRP row 136
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H136").Value = 7362900.5
$ws.Range("I136").Value = 11906703
$ws.Range("J136").Value = 1001577.5
$ws.Range("K136").Value = 35720109
$ws.Range("L136").Value = 3004732.5
$ws.Range("M136").Value = -35717559
$ws.Range("N136").Value = -3009832.5

# CUL row 87
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H87").Value = 6208.4287
$ws.Range("I87").Value = 5762.923
$ws.Range("J87").Value = 12000
$ws.Range("K87").Value = 17288.769
$ws.Range("L87").Value = 36000
$ws.Range("M87").Value = -16040.769
$ws.Range("N87").Value = -38496

# CUL row 90
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H90").Value = 6208.4287
$ws.Range("I90").Value = 5762.923
$ws.Range("J90").Value = 12000
$ws.Range("K90").Value = 51866.307
$ws.Range("L90").Value = 108000
$ws.Range("M90").Value = -45626.307
$ws.Range("N90").Value = -120480

# GSM row 126
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H126").Value = 10329.782
$ws.Range("I126").Value = 13505.625
$ws.Range("J126").Value = 3070.7144
$ws.Range("K126").Value = 40516.875
$ws.Range("L126").Value = 9212.143199999999
$ws.Range("M126").Value = -38046.875
$ws.Range("N126").Value = -14152.1432

# GSM row 133
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H133").Value = 46220
$ws.Range("J133").Value = 46220
$ws.Range("L133").Value = 46220
$ws.Range("N133").Value = -56340

# LTW row 7
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 3774.9092
$ws.Range("I7").Value = 2167.3333
$ws.Range("J7").Value = 5704
$ws.Range("K7").Value = 2167.3333
$ws.Range("L7").Value = 5704
$ws.Range("M7").Value = -2055.3333
$ws.Range("N7").Value = -5928

# LTW row 126
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H126").Value = 3774.9092
$ws.Range("I126").Value = 2167.3333
$ws.Range("J126").Value = 5704
$ws.Range("K126").Value = 6501.999899999999
$ws.Range("L126").Value = 17112
$ws.Range("M126").Value = -4031.999899999999
$ws.Range("N126").Value = -22052

# LTW row 136
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H136").Value = 7790.2256
$ws.Range("I136").Value = 3077.389
$ws.Range("J136").Value = 14315.692
$ws.Range("K136").Value = 9232.167000000001
$ws.Range("L136").Value = 42947.076
$ws.Range("M136").Value = -6682.167000000001
$ws.Range("N136").Value = -48047.076

# LTW row 140
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H140").Value = 80305.14
$ws.Range("J140").Value = 80305.14
$ws.Range("L140").Value = 80305.14
$ws.Range("N140").Value = -90665.14

# WVR row 75
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H75").Value = 41065
$ws.Range("J75").Value = 41065
$ws.Range("L75").Value = 41065
$ws.Range("N75").Value = -42937

# WVR row 78
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H78").Value = 41065
$ws.Range("J78").Value = 41065
$ws.Range("L78").Value = 123195
$ws.Range("N78").Value = -132555

# WVR row 94
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H94").Value = 16866
$ws.Range("I94").Value = 5000
$ws.Range("J94").Value = 19832.5
$ws.Range("K94").Value = 5000
$ws.Range("L94").Value = 19832.5
$ws.Range("M94").Value = -4099
$ws.Range("N94").Value = -21634.5

# WVR row 96
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H96").Value = 3350.375
$ws.Range("I96").Value = 2467.1667
$ws.Range("J96").Value = 6000
$ws.Range("K96").Value = 2467.1667
$ws.Range("L96").Value = 6000
$ws.Range("M96").Value = -1094.1667
$ws.Range("N96").Value = -8746

# WVR row 122
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 974.2143
$ws.Range("I122").Value = 954.8889
$ws.Range("J122").Value = 1009
$ws.Range("K122").Value = 2864.6667
$ws.Range("L122").Value = 3027
$ws.Range("M122").Value = -414.6667000000002
$ws.Range("N122").Value = -7927

# WVR row 126
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 655.6
$ws.Range("I126").Value = 591.3333
$ws.Range("J126").Value = 872.5
$ws.Range("K126").Value = 1773.9999
$ws.Range("L126").Value = 2617.5
$ws.Range("M126").Value = 696.0001
$ws.Range("N126").Value = -7557.5
